$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from column J (9..) into column K for rows 3-6, then set 2023 values
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 975
$ws.Range("K5").Value = 240
$ws.Range("K6").Value = 735
